# Generated PowerShell cell assignments
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing rows 51-54 so sheet shrinks to 50 rows
$ws.Range("A51:D54").EntireRow.Delete()

# New column D header needs the same bold/border/center style as B1:C1 (style index 1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range('B1').Value = 'feat'
$ws.Range('C1').Value = 'shap'
$ws.Range('D1').Value = 'rename'
$ws.Range('A2').Value = 0
$ws.Range('B2').Value = 'From_Same_Orbital'
$ws.Range('C2').Value = 0.00008943528043591508
$ws.Range('D2').Value = '$\mathbf{b}$'
$ws.Range('A3').Value = 2
$ws.Range('B3').Value = 'coulomb'
$ws.Range('C3').Value = 0.0004390300724731704
$ws.Range('D3').Value = '$\langle pp \vert \vert qq \rangle$'
$ws.Range('A4').Value = 3
$ws.Range('B4').Value = 'screen1_1'
$ws.Range('C4').Value = 0.00006731817323400504
$ws.Range('D4').Value = '$(\langle pp \vert \vert rr \rangle)_{1}$'
$ws.Range('A5').Value = 4
$ws.Range('B5').Value = 'screen1_2'
$ws.Range('C5').Value = 0.00008123267920030414
$ws.Range('D5').Value = '$(\langle pp \vert \vert rr \rangle)_{2}$'
$ws.Range('A6').Value = 5
$ws.Range('B6').Value = 'screen1_3'
$ws.Range('C6').Value = 0.0001602502889865089
$ws.Range('D6').Value = '$(\langle pp \vert \vert rr \rangle)_{3}$'
$ws.Range('A7').Value = 6
$ws.Range('B7').Value = 'screen1_4'
$ws.Range('C7').Value = 0.0001656290100168945
$ws.Range('D7').Value = '$(\langle pp \vert \vert rr \rangle)_{4}$'
$ws.Range('A8').Value = 7
$ws.Range('B8').Value = 'screen2_1'
$ws.Range('C8').Value = 0.0001471514837970222
$ws.Range('D8').Value = '$(\langle qq \vert \vert ss \rangle)_{1}$'
$ws.Range('A9').Value = 8
$ws.Range('B9').Value = 'screen2_2'
$ws.Range('C9').Value = 0.0001438513075941875
$ws.Range('D9').Value = '$(\langle qq \vert \vert ss \rangle)_{2}$'
$ws.Range('A10').Value = 9
$ws.Range('B10').Value = 'screen2_3'
$ws.Range('C10').Value = 0.003041246658121699
$ws.Range('D10').Value = '$(\langle qq \vert \vert ss \rangle)_{3}$'
$ws.Range('A11').Value = 10
$ws.Range('B11').Value = 'screen2_4'
$ws.Range('C11').Value = 0.00009364611295209535
$ws.Range('D11').Value = '$(\langle qq \vert \vert ss \rangle)_{4}$'
$ws.Range('A12').Value = 12
$ws.Range('B12').Value = 'eijab_2'
$ws.Range('C12').Value = 0.00001626226690895013
$ws.Range('D12').Value = '$(e^{rs}_{pq})_{2}$'
$ws.Range('A13').Value = 13
$ws.Range('B13').Value = 'eijab_3'
$ws.Range('C13').Value = 0.00002192378398809872
$ws.Range('D13').Value = '$(e^{rs}_{pq})_{3}$'
$ws.Range('A14').Value = 14
$ws.Range('B14').Value = 'eijab_4'
$ws.Range('C14').Value = 0.00001793164792377623
$ws.Range('D14').Value = '$(e^{rs}_{pq})_{4}$'
$ws.Range('A15').Value = 15
$ws.Range('B15').Value = 'screenvirt_1'
$ws.Range('C15').Value = 0.00007753486198487605
$ws.Range('D15').Value = '$(\langle ss \vert \vert rr \rangle)_{1}$'
$ws.Range('A16').Value = 16
$ws.Range('B16').Value = 'screenvirt_2'
$ws.Range('C16').Value = 0.00008887849039566642
$ws.Range('D16').Value = '$(\langle ss \vert \vert rr \rangle)_{2}$'
$ws.Range('A17').Value = 17
$ws.Range('B17').Value = 'screenvirt_3'
$ws.Range('C17').Value = 0.0001935646786050618
$ws.Range('D17').Value = '$(\langle ss \vert \vert rr \rangle)_{3}$'
$ws.Range('A18').Value = 18
$ws.Range('B18').Value = 'screenvirt_4'
$ws.Range('C18').Value = 0.00009128295436749971
$ws.Range('D18').Value = '$(\langle ss \vert \vert rr \rangle)_{4}$'
$ws.Range('A19').Value = 19
$ws.Range('B19').Value = 'Fr1'
$ws.Range('C19').Value = 0.001615012776455012
$ws.Range('D19').Value = '$(F_{r})_{1}$'
$ws.Range('A20').Value = 20
$ws.Range('B20').Value = 'Fr2'
$ws.Range('C20').Value = 0.0002321659468906043
$ws.Range('D20').Value = '$(F_{r})_{2}$'
$ws.Range('A21').Value = 21
$ws.Range('B21').Value = 'Fr3'
$ws.Range('C21').Value = 0.0001756649717269622
$ws.Range('D21').Value = '$(F_{r})_{3}$'
$ws.Range('A22').Value = 22
$ws.Range('B22').Value = 'Fr4'
$ws.Range('C22').Value = 0.0003483228146047222
$ws.Range('D22').Value = '$(F_{r})_{4}$'
$ws.Range('A23').Value = 23
$ws.Range('B23').Value = 'Fs1'
$ws.Range('C23').Value = 0.0000621358617630117
$ws.Range('D23').Value = '$(F_{s})_{1}$'
$ws.Range('A24').Value = 24
$ws.Range('B24').Value = 'Fs2'
$ws.Range('C24').Value = 0.00001106217610872746
$ws.Range('D24').Value = '$(F_{s})_{2}$'
$ws.Range('A25').Value = 27
$ws.Range('B25').Value = 'occr1'
$ws.Range('C25').Value = 0.000066433079736839
$ws.Range('D25').Value = '$(\eta_{r})_{1}$'
$ws.Range('A26').Value = 28
$ws.Range('B26').Value = 'occr2'
$ws.Range('C26').Value = 0.00003687779246704264
$ws.Range('D26').Value = '$(\eta_{r})_{2}$'
$ws.Range('A27').Value = 30
$ws.Range('B27').Value = 'occr4'
$ws.Range('C27').Value = 0.00004369366155118692
$ws.Range('D27').Value = '$(\eta_{r})_{4}$'
$ws.Range('A28').Value = 35
$ws.Range('B28').Value = 'SCFFr1'
$ws.Range('C28').Value = 0.0001109908010416703
$ws.Range('D28').Value = '$(F_{r}^{\text{SCF}})_{1}$'
$ws.Range('A29').Value = 36
$ws.Range('B29').Value = 'SCFFr2'
$ws.Range('C29').Value = 0.0002771385770950991
$ws.Range('D29').Value = '$(F_{r}^{\text{SCF}})_{2}$'
$ws.Range('A30').Value = 37
$ws.Range('B30').Value = 'SCFFr3'
$ws.Range('C30').Value = 0.00005507133313061352
$ws.Range('D30').Value = '$(F_{r}^{\text{SCF}})_{3}$'
$ws.Range('A31').Value = 38
$ws.Range('B31').Value = 'SCFFr4'
$ws.Range('C31').Value = 0.00004177197385700876
$ws.Range('D31').Value = '$(F_{r}^{\text{SCF}})_{4}$'
$ws.Range('A32').Value = 39
$ws.Range('B32').Value = 'SCFFs1'
$ws.Range('C32').Value = 0.00002694116138362471
$ws.Range('D32').Value = '$(F_{s}^{\text{SCF}})_{1}$'
$ws.Range('A33').Value = 40
$ws.Range('B33').Value = 'SCFFs2'
$ws.Range('C33').Value = 0.00004164160574885197
$ws.Range('D33').Value = '$(F_{s}^{\text{SCF}})_{2}$'
$ws.Range('A34').Value = 41
$ws.Range('B34').Value = 'SCFFs3'
$ws.Range('C34').Value = 0.00005302720329996195
$ws.Range('D34').Value = '$(F_{s}^{\text{SCF}})_{3}$'
$ws.Range('A35').Value = 42
$ws.Range('B35').Value = 'SCFFs4'
$ws.Range('C35').Value = 0.00003854034152464111
$ws.Range('D35').Value = '$(F_{s}^{\text{SCF}})_{4}$'
$ws.Range('A36').Value = 51
$ws.Range('B36').Value = 'hrr1'
$ws.Range('C36').Value = 0.00007587186744523922
$ws.Range('D36').Value = '$(h_{rr})_{1}$'
$ws.Range('A37').Value = 52
$ws.Range('B37').Value = 'hrr2'
$ws.Range('C37').Value = 0.00003655237044401989
$ws.Range('D37').Value = '$(h_{rr})_{2}$'
$ws.Range('A38').Value = 53
$ws.Range('B38').Value = 'hrr3'
$ws.Range('C38').Value = 0.0001062744063263749
$ws.Range('D38').Value = '$(h_{rr})_{3}$'
$ws.Range('A39').Value = 54
$ws.Range('B39').Value = 'hrr4'
$ws.Range('C39').Value = 0.00007542162687745788
$ws.Range('D39').Value = '$(h_{rr})_{4}$'
$ws.Range('A40').Value = 55
$ws.Range('B40').Value = 'hss1'
$ws.Range('C40').Value = 0.00006136511770506882
$ws.Range('D40').Value = '$(h_{ss})_{1}$'
$ws.Range('A41').Value = 56
$ws.Range('B41').Value = 'hss2'
$ws.Range('C41').Value = 0.00009783626529416461
$ws.Range('D41').Value = '$(h_{ss})_{2}$'
$ws.Range('A42').Value = 57
$ws.Range('B42').Value = 'hss3'
$ws.Range('C42').Value = 0.00002475638289042079
$ws.Range('D42').Value = '$(h_{ss})_{3}$'
$ws.Range('A43').Value = 58
$ws.Range('B43').Value = 'hss4'
$ws.Range('C43').Value = 0.00004186632922329535
$ws.Range('D43').Value = '$(h_{ss})_{4}$'
$ws.Range('A44').Value = 59
$ws.Range('B44').Value = 'hpp'
$ws.Range('C44').Value = 0.0008505911323295213
$ws.Range('D44').Value = '$h_{pp}$'
$ws.Range('A45').Value = 61
$ws.Range('B45').Value = 'Fp'
$ws.Range('C45').Value = 0.001759077041079878
$ws.Range('D45').Value = '$F_{p}$'
$ws.Range('A46').Value = 62
$ws.Range('B46').Value = 'Fq'
$ws.Range('C46').Value = 0.0005954976112079823
$ws.Range('D46').Value = '$F_{q}$'
$ws.Range('A47').Value = 63
$ws.Range('B47').Value = 'occp'
$ws.Range('C47').Value = 0.00001795169140818669
$ws.Range('D47').Value = '$\eta_{p}$'
$ws.Range('A48').Value = 64
$ws.Range('B48').Value = 'occq'
$ws.Range('C48').Value = 0.00004166667602412304
$ws.Range('D48').Value = '$\eta_{q}$'
$ws.Range('A49').Value = 65
$ws.Range('B49').Value = 'SCFFp'
$ws.Range('C49').Value = 0.00007323152118066662
$ws.Range('D49').Value = '$F_{p}^{\text{SCF}}$'
$ws.Range('A50').Value = 66
$ws.Range('B50').Value = 'SCFFq'
$ws.Range('C50').Value = 0.00103208798389803
$ws.Range('D50').Value = '$F_{q}^{\text{SCF}}$'

Write-Output "applied edits"
